$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

$ws.Range("A6").Value = "Amor"
$ws.Range("B6").Value = "amor9090"
